$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3.82101041052307,9.913587179101524,13.92254190664093,28.90802875049724,27.85331984228921,14.07418700545614,19.87782691889145,9.923716077852141,11.96212016468163,15.95276652822873,17.43466329006858,21.30072375762587),
    @(3.688339742361125,9.898883860528752,13.94012509482078,28.95805630891429,27.88608789775542,14.11821499009167,19.92072875713331,9.949626404663608,11.53576468077759,15.78908034129755,17.46298577517745,21.36531929834431),
    @(3.603562927336873,9.891349186623698,13.95327599606256,28.99591390959203,27.91566989850799,14.14756801349818,19.95176989224539,9.966671902403034,11.26638721992872,15.68969374974176,17.48204306711567,21.40975639174475),
    @(3.568211111053967,9.888656838473596,13.95922735483805,29.01313338528085,27.93009699873783,14.16011268366748,19.96559873357975,9.973904191213665,11.1548596308884,15.64951202926436,17.49022902196183,21.42906297725789),
    @(3.562293299222381,9.888232683709351,13.96025134601605,29.01610080106174,27.93263564253213,14.1622309262245,19.96796615657291,9.975122397650367,11.13623946445608,15.64286026348694,17.49161367911521,21.43234110379517),
    @(3.603089376796725,9.891311342510226,13.95335386003742,28.99613888578929,27.91585487313878,14.14773483476356,19.95195162109465,9.966768280675968,11.26489000532817,15.6891505029043,17.48215176442328,21.41001191961954),
    @(3.775970502048753,9.908209188073615,13.92811597893193,28.92379441702926,27.86265052381978,14.08888627790107,19.89164283009806,9.932414292408092,11.8167772142931,15.89611864886713,17.44408329904511,21.32200366454724),
    @(4.08754958741915,9.953070362616453,13.89729859500406,28.83870850896125,27.83364533407634,13.99189964662943,19.81075657838694,9.874048198576201,12.83273003723835,16.30915424418472,17.38262767752385,21.18742161802498),
    @(4.298444482221182,9.992999104452256,13.8860221551826,28.81094572096221,27.8584920515948,13.93188285929636,19.77422618806691,9.836633453480742,13.53167414155224,16.61470260100158,17.34547897141247,21.11185764961239),
    @(4.390253657676321,10.01263438776058,13.88335319001058,28.80587611663303,27.87983182643482,13.90702280953067,19.76259588805151,9.820795192769992,13.83811478443656,16.75367551656691,17.33030787023764,21.08257195539633),
    @(4.424410625038073,10.02027736911547,13.88269551753668,28.80504342302913,27.88935404332814,13.89796032927901,19.75890971583152,9.814967288828862,13.95241531609756,16.80625768956568,17.32481073096079,21.07221597997951),
    @(4.417081629804922,10.01862215371468,13.88282147411327,28.80517442007432,27.88723921244778,13.89989645906381,19.75967165899782,9.816214888683687,13.92787738089165,16.79493576774437,17.32598362661814,21.07441365330299),
    @(4.393076062419738,10.01325904540822,13.88329201351168,28.80578582712566,27.88058635699524,13.90627018772241,19.76227823224749,9.820312327968356,13.84755372357207,16.75800259274785,17.329850654244,21.08170524352877),
    @(4.378292200617486,10.0100008881158,13.88362617664624,28.80630188427673,27.87669890099189,13.91022006363863,19.76396835063223,9.822844220554822,13.79812394413508,16.73537305144695,17.33225157511668,21.08626718577091),
    @(4.292360090656187,9.991745159497405,13.88624602606723,28.8114291685159,27.85729928776682,13.93355668847082,19.77508669214317,9.837692282059452,13.51140822090114,16.60561645436865,17.34650515447267,21.11387412344347),
    @(4.238574479018684,9.980919917591065,13.88848292197714,28.81651105605722,27.84796811865266,13.94849862469639,19.78318553406851,9.847103616891344,13.33250369914085,16.52598011949366,17.35569138406479,21.13211501020133),
    @(4.207250660410045,9.974832286904423,13.8900012025277,28.82014563980877,27.84354603213076,13.9573226537136,19.78831320461218,9.852628028163316,13.22852334924063,16.4801756777597,17.36113775049613,21.1430855468262),
    @(4.196578866012269,9.972795072866361,13.89055508023677,28.82149845707488,27.84221111652078,13.96034977996182,19.7901299413585,9.854517615763204,13.1931348532272,16.46466837495623,17.36300975890413,21.14688215874783),
    @(4.244340300245135,9.982057950542478,13.88822082850587,28.8158964286967,27.84886364666187,13.9468842417937,19.78227480741447,9.84609025104603,13.35166074545102,16.5344578113303,17.35469666120589,21.13012366028376),
    @(4.400143725953599,10.01482872166276,13.88314423168829,28.80557674373521,27.88250137029771,13.90438852737514,19.76149312881866,9.819104207072538,13.87119463686322,16.76885227206982,17.32870809380544,21.07954359345435),
    @(4.498413757772483,10.03745389458067,13.88188347364028,28.80516787030503,27.91288437093014,13.87866396976382,19.75209605922799,9.802456303642582,14.20055452761983,16.92176864989127,17.31316728259362,21.05076505938273),
    @(4.44629522272847,10.02526932116793,13.88236845927866,28.80480659121179,27.89590103812813,13.89220607174597,19.75672837562027,9.811251188817293,14.02572666262286,16.84019260575626,17.32132977881135,21.06573257518261),
    @(4.241734822650385,9.981543022093396,13.8883385974447,28.81617208098329,27.84845584242412,13.94761337633012,19.78268507818964,9.846548039530036,13.34300334461443,16.53062510744387,17.35514586114284,21.13102244333454),
    @(4.006341365750074,9.939696524130152,13.90363773805988,28.85563109131682,27.83339651704417,14.01616425657591,19.82862411539255,9.888876239267125,12.56576901347881,16.19688487238862,17.38262767752385,21.21974546932175)
)

$cols = @(3,4,5,6,7,8,9,10,11,13,14,15)
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25)

for ($r = 0; $r -lt $rows.Length; $r++) {
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $ws.Cells.Item($rows[$r], $cols[$c]).Value = $data[$r][$c]
    }
}

Write-Output "applied"